$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# PageElement sheet: add the XPATH pairs needed for the new
# "ClueTransferCustomer" test (transfer a clue into a customer, then
# fill out + save the "Add customer" form).
# ---------------------------------------------------------------------
$pe = $wb.Worksheets.Item("PageElement")

$pe.Range("A16").Value = "ViewCluePage_TransferClue1"
$pe.Range("B16").Value = "//span[text()='"

$pe.Range("A17").Value = "ViewCluePage_TransferClue2"
# Leading "'" is Excel's quote-prefix marker and is swallowed by Value;
# double it so the literal apostrophe that starts the XPath survives.
$pe.Range("B17").Value = "'']/../../..//a[contains(text(),'转换')]"

$pe.Range("A18").Value = "AddCustomerPage_Name"
$pe.Range("B18").Value = "//input[@id='name']"

$pe.Range("A19").Value = "AddCustomerPage_Industry"
$pe.Range("B19").Value = "//input[@id='industry']"

$pe.Range("A20").Value = "AddCustomerPage_ComOrig"
$pe.Range("B20").Value = "//select[@id='origin']"

$pe.Range("A21").Value = "AddCustomerPage_ComOrigSelection"
$pe.Range("B21").Value = "//select[@id='origin']/option[@value='网络营销']"

$pe.Range("A22").Value = "AddCustomerPage_ComPro"
$pe.Range("B22").Value = "//input[@id='ownership2']"

$pe.Range("A23").Value = "AddCustomerPage_ConName"
$pe.Range("B23").Value = "//input[@name='con_name']"

$pe.Range("A24").Value = "AddCustomerPage_NumEmp"
$pe.Range("B24").Value = "//select[@id='no_of_employees']"

$pe.Range("A25").Value = "AddCustomerPage_NumEmpSelection"
$pe.Range("B25").Value = "//select[@id='no_of_employees']/option[@value='5--20人']"

$pe.Range("A26").Value = "AddCustomerPage_CreateBussiness"
$pe.Range("B26").Value = "//input[@name='create_business2']"

$pe.Range("A27").Value = "AddCustomerPage_SaveBtn"
$pe.Range("B27").Value = "//input[@value='保存']"

$pe.Range("A16:B27").Select() | Out-Null
$pe.Range("B27").Select() | Out-Null

# ---------------------------------------------------------------------
# TestData sheet: verification text for the "add customer" step, and
# the expected customer name/value pair used by the new script.
# ---------------------------------------------------------------------
$td = $wb.Worksheets.Item("TestData")

$td.Range("D9").Value = "VerifyAddCustomerPage"
$td.Range("E9").Value = "添加客户"

$td.Range("A10").Value = "CustomerName"
$td.Range("B10").Value = "MaoTaoCEO"
$td.Range("D10").Value = "VerifyAddCustomerResult"
$td.Range("E10").Value = "添加客户成功"

$td.Range("E10").Select() | Out-Null
